$d = $word.ActiveDocument

# --- Hunk 1 -----------------------------------------------------------
# The runs "<id>" + "p139r_1" + "</id>" (three separately-formatted runs)
# collapse into a single run "<id>p139r_1</id>" using the formatting of
# the first run (Courier New / color 7f6000 / sz 18). A Find/Replace of
# the full matched text with itself causes Word to re-run/merge the
# matched range into one run using the first run's character formatting.
$rng1 = $d.Content
$rng1.Find.Execute("<id>p139r_1</id>", $false, $false, $false, $false, $false, `
                    $true, 1, $false, "<id>p139r_1</id>", 2) | Out-Null

# --- Hunk 2 -----------------------------------------------------------
# "ligne de " (typo) becomes "ligue de " (correction), with the result
# represented as three runs: "lig" (unchanged formatting), "u" (the
# corrected letter, losing the explicit black color) and "e de "
# (unchanged formatting). We locate the run, fix the single letter, and
# then apply distinguishing character formatting so the single changed
# letter splits into its own run.
$rng2 = $d.Content
$rng2.Find.Execute("ligne de ", $false, $false, $false, $false, $false, `
                    $true, 1, $false, "", 0) | Out-Null

$start2 = $rng2.Start
$letter = $d.Range($start2 + 3, $start2 + 4)
$letter.Text = "u"

# Re-fetch the (now retargeted) one-character range for the replaced
# letter and give it its own (automatic-color) run so it splits away
# from the unchanged "lig" / "e de " runs around it.
$letter2 = $d.Range($start2 + 3, $start2 + 4)
$letter2.Font.Color = -16777216
